# Update the repayment strategy selection on the ProductLoanInput sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Activate()
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"
$ws.Range("B17").Select()
